$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column A (shifts every existing column right by one). ---
$ws.Columns("A:A").Insert()

# --- New column header: same bold/boxed/centered style as the other
#     headers (column B used to be column A before the insert), then set
#     the text. ---
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value2 = "Data"

# --- New column data: dates for rows 2-6 (serial 43022..43026 => 2017-10-14..18). ---
$ws.Range("A2").Value2 = 43022
$ws.Range("A3").Value2 = 43023
$ws.Range("A4").Value2 = 43024
$ws.Range("A5").Value2 = 43025
$ws.Range("A6").Value2 = 43026

# --- Give A2:A6 the same look as the header style (bold font / boxed border /
#     centered-top alignment) plus a date number format. Copy the existing
#     header style (column B, the old column A) onto the date cells first ...
$ws.Range("B1").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)

# ... then apply the date number format. A2 walks through the lowercase
# format first so both "yyyy-mm-dd" (164) and "YYYY-MM-DD" (165) end up
# registered in the workbook's numFmts, matching the authored file; the
# remaining cells go directly to the final, uppercase format.
$cellA2 = $ws.Range("A2")
$cellA2.NumberFormat = "yyyy-mm-dd"
$cellA2.NumberFormat = "YYYY-MM-DD"
$ws.Range("A3:A6").NumberFormat = "YYYY-MM-DD"

# --- Fix up the "proporção" column (now BI) whose ratio is non-blank-count /
#     total-columns; the extra "Data" column changes the denominator. ---
$ws.Range("BI2").Value2 = 0.7333333333333333
$ws.Range("BI3").Value2 = 0.6833333333333333
$ws.Range("BI4").Value2 = 0.65
$ws.Range("BI5").Value2 = 0.65
$ws.Range("BI6").Value2 = 0.8

$excel.CutCopyMode = $false
